$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.62515303497554
$ws.Cells.Item(2, 3).Value = 0.1858221518902496
$ws.Cells.Item(2, 4).Value = 0.1330609785644725
$ws.Cells.Item(2, 5).Value = 0.1137359540516814
$ws.Cells.Item(2, 6).Value = 1.471928477265322
$ws.Cells.Item(2, 8).Value = 0.07973214163530429
$ws.Cells.Item(2, 9).Value = 0.8647057324364056
$ws.Cells.Item(2, 10).Value = 0.1314700272971461
$ws.Cells.Item(2, 12).Value = 0.3616644208880757
$ws.Cells.Item(2, 15).Value = 3.715882654557362

$ws.Cells.Item(3, 2).Value = 1.48849865732393
$ws.Cells.Item(3, 3).Value = 0.1663086030168017
$ws.Cells.Item(3, 4).Value = 0.1320441260592986
$ws.Cells.Item(3, 5).Value = 0.1145093657757812
$ws.Cells.Item(3, 6).Value = 1.483164092666136
$ws.Cells.Item(3, 8).Value = 0.07973214163530429
$ws.Cells.Item(3, 9).Value = 0.8782526259561507
$ws.Cells.Item(3, 10).Value = 0.1329526615391807
$ws.Cells.Item(3, 12).Value = 0.3517483801029186
$ws.Cells.Item(3, 15).Value = 3.751506263169773

$ws.Cells.Item(4, 2).Value = 1.404576991406657
$ws.Cells.Item(4, 3).Value = 0.1542672123198372
$ws.Cells.Item(4, 4).Value = 0.131452016654066
$ws.Cells.Item(4, 5).Value = 0.1150229574208899
$ws.Cells.Item(4, 6).Value = 1.490983473710244
$ws.Cells.Item(4, 8).Value = 0.07973214163530429
$ws.Cells.Item(4, 9).Value = 0.8871437396527604
$ws.Cells.Item(4, 10).Value = 0.133913938860216
$ws.Cells.Item(4, 12).Value = 0.345762628437825
$ws.Cells.Item(4, 15).Value = 3.775907542300246

$ws.Cells.Item(5, 2).Value = 1.370376729839222
$ws.Cells.Item(5, 3).Value = 0.1493454704450699
$ws.Cells.Item(5, 4).Value = 0.1312188800539928
$ws.Cells.Item(5, 5).Value = 0.1152420000337138
$ws.Cells.Item(5, 6).Value = 1.494401347214598
$ws.Cells.Item(5, 8).Value = 0.07973214163530429
$ws.Cells.Item(5, 9).Value = 0.8909109525795458
$ws.Cells.Item(5, 10).Value = 0.1343184849281522
$ws.Cells.Item(5, 12).Value = 0.3433494274882349
$ws.Cells.Item(5, 15).Value = 3.78648643188248

$ws.Cells.Item(6, 2).Value = 1.364697781428106
$ws.Cells.Item(6, 3).Value = 0.1485273355704351
$ws.Cells.Item(6, 4).Value = 0.1311806616425528
$ws.Cells.Item(6, 5).Value = 0.1152789611994285
$ws.Cells.Item(6, 6).Value = 1.494982856386834
$ws.Cells.Item(6, 8).Value = 0.07973214163530429
$ws.Cells.Item(6, 9).Value = 0.8915451892553783
$ws.Cells.Item(6, 10).Value = 0.1343864339507235
$ws.Cells.Item(6, 12).Value = 0.3429502962744664
$ws.Cells.Item(6, 15).Value = 3.788281394248543

$ws.Cells.Item(7, 2).Value = 1.404115757664783
$ws.Cells.Item(7, 3).Value = 0.154200895409474
$ws.Cells.Item(7, 4).Value = 0.1314488394216369
$ws.Cells.Item(7, 5).Value = 0.1150258720093955
$ws.Cells.Item(7, 6).Value = 1.491028631499567
$ws.Cells.Item(7, 8).Value = 0.07973214163530429
$ws.Cells.Item(7, 9).Value = 0.8871939626932743
$ws.Cells.Item(7, 10).Value = 0.1339193427981489
$ws.Cells.Item(7, 12).Value = 0.3457299774483005
$ws.Cells.Item(7, 15).Value = 3.776047641897321

$ws.Cells.Item(8, 2).Value = 1.578039149006656
$ws.Cells.Item(8, 3).Value = 0.1791065100538276
$ws.Cells.Item(8, 4).Value = 0.1327037025413276
$ws.Cells.Item(8, 5).Value = 0.1139946003789012
$ws.Cells.Item(8, 6).Value = 1.475611373455699
$ws.Cells.Item(8, 8).Value = 0.07973214163530429
$ws.Cells.Item(8, 9).Value = 0.8692576367297775
$ws.Cells.Item(8, 10).Value = 0.1319706744064132
$ws.Cells.Item(8, 12).Value = 0.3582241672246624
$ws.Cells.Item(8, 15).Value = 3.727640591528683

$ws.Cells.Item(9, 2).Value = 1.918894934819832
$ws.Cells.Item(9, 3).Value = 0.2274591484837742
$ws.Cells.Item(9, 4).Value = 0.1354185305523501
$ws.Cells.Item(9, 5).Value = 0.1122787943206909
$ws.Cells.Item(9, 6).Value = 1.452688390071152
$ws.Cells.Item(9, 8).Value = 0.07973214163530429
$ws.Cells.Item(9, 9).Value = 0.838639107783866
$ws.Cells.Item(9, 10).Value = 0.1285529823016447
$ws.Cells.Item(9, 12).Value = 0.3835328042511748
$ws.Cells.Item(9, 15).Value = 3.652802128671567

$ws.Cells.Item(10, 2).Value = 2.169105577986045
$ws.Cells.Item(10, 3).Value = 0.2626752770323151
$ws.Cells.Item(10, 4).Value = 0.1375658956395753
$ws.Cells.Item(10, 5).Value = 0.1112041559391592
$ws.Cells.Item(10, 6).Value = 1.440311556566556
$ws.Cells.Item(10, 8).Value = 0.07973214163530429
$ws.Cells.Item(10, 9).Value = 0.8189279029201124
$ws.Cells.Item(10, 10).Value = 0.1262873216944915
$ws.Cells.Item(10, 12).Value = 0.4026113619403588
$ws.Cells.Item(10, 15).Value = 3.610103593473184

$ws.Cells.Item(11, 2).Value = 2.282868830489917
$ws.Cells.Item(11, 3).Value = 0.2786268128057827
$ws.Cells.Item(11, 4).Value = 0.1385755415547436
$ws.Cells.Item(11, 5).Value = 0.110755466858242
$ws.Cells.Item(11, 6).Value = 1.435652315916258
$ws.Cells.Item(11, 8).Value = 0.07973214163530429
$ws.Cells.Item(11, 9).Value = 0.8105669614605233
$ws.Cells.Item(11, 10).Value = 0.1253097210512319
$ws.Cells.Item(11, 12).Value = 0.4113943097949715
$ws.Cells.Item(11, 15).Value = 3.593355146968548

$ws.Cells.Item(12, 2).Value = 2.325937589928742
$ws.Cells.Item(12, 3).Value = 0.2846571441912431
$ws.Cells.Item(12, 4).Value = 0.1389625437587227
$ws.Cells.Item(12, 5).Value = 0.1105913215113077
$ws.Cells.Item(12, 6).Value = 1.434027754013186
$ws.Cells.Item(12, 8).Value = 0.07973214163530429
$ws.Cells.Item(12, 9).Value = 0.8074881550827833
$ws.Cells.Item(12, 10).Value = 0.1249471488524125
$ws.Cells.Item(12, 12).Value = 0.4147349662544428
$ws.Cells.Item(12, 15).Value = 3.587398342472483

$ws.Cells.Item(13, 2).Value = 2.31666248632223
$ws.Cells.Item(13, 3).Value = 0.2833588618400711
$ws.Cells.Item(13, 4).Value = 0.1388789888743887
$ws.Cells.Item(13, 5).Value = 0.1106264170206952
$ws.Cells.Item(13, 6).Value = 1.434371412743658
$ws.Cells.Item(13, 8).Value = 0.07973214163530429
$ws.Cells.Item(13, 9).Value = 0.8081473453840466
$ws.Cells.Item(13, 10).Value = 0.1250248962797762
$ws.Cells.Item(13, 12).Value = 0.4140148437109872
$ws.Cells.Item(13, 15).Value = 3.588664091021258

$ws.Cells.Item(14, 2).Value = 2.286412356000824
$ws.Cells.Item(14, 3).Value = 0.2791231374507959
$ws.Cells.Item(14, 4).Value = 0.1386072870991129
$ws.Cells.Item(14, 5).Value = 0.1107418470865351
$ws.Cells.Item(14, 6).Value = 1.435515859709767
$ws.Cells.Item(14, 8).Value = 0.07973214163530429
$ws.Cells.Item(14, 9).Value = 0.8103119154856238
$ws.Cells.Item(14, 10).Value = 0.1252797392568534
$ws.Cells.Item(14, 12).Value = 0.4116688533201369
$ws.Cells.Item(14, 15).Value = 3.59285734737341

$ws.Cells.Item(15, 2).Value = 2.267881783702705
$ws.Cells.Item(15, 3).Value = 0.2765273000628952
$ws.Cells.Item(15, 4).Value = 0.1384414689088516
$ws.Cells.Item(15, 5).Value = 0.1108133015288182
$ws.Cells.Item(15, 6).Value = 1.436235076118408
$ws.Cells.Item(15, 8).Value = 0.07973214163530429
$ws.Cells.Item(15, 9).Value = 0.8116491519248825
$ws.Cells.Item(15, 10).Value = 0.1254368306267748
$ws.Cells.Item(15, 12).Value = 0.4102337803747531
$ws.Cells.Item(15, 15).Value = 3.59547605984838

$ws.Cells.Item(16, 2).Value = 2.161669487703932
$ws.Cells.Item(16, 3).Value = 0.2616314012662428
$ws.Cells.Item(16, 4).Value = 0.1375005690288518
$ws.Cells.Item(16, 5).Value = 0.111234286037357
$ws.Cells.Item(16, 6).Value = 1.440635602181956
$ws.Cells.Item(16, 8).Value = 0.07973214163530429
$ws.Cells.Item(16, 9).Value = 0.819486516255381
$ws.Cells.Item(16, 10).Value = 0.1263522773797019
$ws.Cells.Item(16, 12).Value = 0.4020394505036222
$ws.Cells.Item(16, 15).Value = 3.611252043332968

$ws.Cells.Item(17, 2).Value = 2.096494863185853
$ws.Cells.Item(17, 3).Value = 0.2524754839970456
$ws.Cells.Item(17, 4).Value = 0.1369317257235849
$ws.Cells.Item(17, 5).Value = 0.1115028255791959
$ws.Cells.Item(17, 6).Value = 1.443584001595212
$ws.Cells.Item(17, 8).Value = 0.07973214163530429
$ws.Cells.Item(17, 9).Value = 0.8244497763617211
$ws.Cells.Item(17, 10).Value = 0.1269274598907597
$ws.Cells.Item(17, 12).Value = 0.3970389838449933
$ws.Cells.Item(17, 15).Value = 3.621615856110992

$ws.Cells.Item(18, 2).Value = 2.059002751585695
$ws.Cells.Item(18, 3).Value = 0.2472028182027373
$ws.Cells.Item(18, 4).Value = 0.1366076319622422
$ws.Cells.Item(18, 5).Value = 0.1116610641107858
$ws.Cells.Item(18, 6).Value = 1.445371232821472
$ws.Cells.Item(18, 8).Value = 0.07973214163530429
$ws.Cells.Item(18, 9).Value = 0.8273615202789379
$ws.Cells.Item(18, 10).Value = 0.1272632842406973
$ws.Cells.Item(18, 12).Value = 0.3941726508161878
$ws.Cells.Item(18, 15).Value = 3.627828616754186

$ws.Cells.Item(19, 2).Value = 2.046307716041781
$ws.Cells.Item(19, 3).Value = 0.2454164907250345
$ws.Cells.Item(19, 4).Value = 0.1364984314523809
$ws.Cells.Item(19, 5).Value = 0.1117152908591805
$ws.Cells.Item(19, 6).Value = 1.445992050116431
$ws.Cells.Item(19, 8).Value = 0.07973214163530429
$ws.Cells.Item(19, 9).Value = 0.8283571721349787
$ws.Cells.Item(19, 10).Value = 0.1273778466038582
$ws.Cells.Item(19, 12).Value = 0.3932038499031307
$ws.Cells.Item(19, 15).Value = 3.629975369257721

$ws.Cells.Item(20, 2).Value = 2.103433389079896
$ws.Cells.Item(20, 3).Value = 0.2534508142800576
$ws.Cells.Item(20, 4).Value = 0.1369919606047674
$ws.Cells.Item(20, 5).Value = 0.1114738477890249
$ws.Cells.Item(20, 6).Value = 1.443260679795756
$ws.Cells.Item(20, 8).Value = 0.07973214163530429
$ws.Cells.Item(20, 9).Value = 0.8239155273623027
$ws.Cells.Item(20, 10).Value = 0.1268657138533538
$ws.Cells.Item(20, 12).Value = 0.3975702790150564
$ws.Cells.Item(20, 15).Value = 3.620486547815972

$ws.Cells.Item(21, 2).Value = 2.295297868255545
$ws.Cells.Item(21, 3).Value = 0.2803675509376831
$ws.Cells.Item(21, 4).Value = 0.1386869660801082
$ws.Cells.Item(21, 5).Value = 0.1107077861611803
$ws.Cells.Item(21, 6).Value = 1.435175912752655
$ws.Cells.Item(21, 8).Value = 0.07973214163530429
$ws.Cells.Item(21, 9).Value = 0.8096737579509146
$ws.Cells.Item(21, 10).Value = 0.1252046788323584
$ws.Cells.Item(21, 12).Value = 0.4123575290863926
$ws.Cells.Item(21, 15).Value = 3.591615218791958

$ws.Cells.Item(22, 2).Value = 2.420627649643507
$ws.Cells.Item(22, 3).Value = 0.2978997729861135
$ws.Cells.Item(22, 4).Value = 0.1398219567277863
$ws.Cells.Item(22, 5).Value = 0.110240708599143
$ws.Cells.Item(22, 6).Value = 1.430706933339806
$ws.Cells.Item(22, 8).Value = 0.07973214163530429
$ws.Cells.Item(22, 9).Value = 0.8008748814609454
$ws.Cells.Item(22, 10).Value = 0.1241635319386583
$ws.Cells.Item(22, 12).Value = 0.422107697588757
$ws.Cells.Item(22, 15).Value = 3.57499329291673

$ws.Cells.Item(23, 2).Value = 2.353743516891029
$ws.Cells.Item(23, 3).Value = 0.2885480397065123
$ws.Cells.Item(23, 4).Value = 0.1392137163239227
$ws.Cells.Item(23, 5).Value = 0.11048692760634
$ws.Cells.Item(23, 6).Value = 1.433017498439753
$ws.Cells.Item(23, 8).Value = 0.07973214163530429
$ws.Cells.Item(23, 9).Value = 0.8055243765186511
$ws.Cells.Item(23, 10).Value = 0.1247151477217905
$ws.Cells.Item(23, 12).Value = 0.4168960641976156
$ws.Cells.Item(23, 15).Value = 3.58365886258872

$ws.Cells.Item(24, 2).Value = 2.100296550960934
$ws.Cells.Item(24, 3).Value = 0.2530098948702175
$ws.Cells.Item(24, 4).Value = 0.1369647192487093
$ws.Cells.Item(24, 5).Value = 0.1114869366499089
$ws.Cells.Item(24, 6).Value = 1.44340656652723
$ws.Cells.Item(24, 8).Value = 0.07973214163530429
$ws.Cells.Item(24, 9).Value = 0.8241568797472638
$ws.Cells.Item(24, 10).Value = 0.1268936132127096
$ws.Cells.Item(24, 12).Value = 0.3973300539720412
$ws.Cells.Item(24, 15).Value = 3.620996315496939

$ws.Cells.Item(25, 2).Value = 1.826715976498406
$ws.Cells.Item(25, 3).Value = 0.2144318424828384
$ws.Cells.Item(25, 4).Value = 0.1346571176349443
$ws.Cells.Item(25, 5).Value = 0.1127102435519678
$ws.Cells.Item(25, 6).Value = 1.458106206644224
$ws.Cells.Item(25, 8).Value = 0.07973214163530429
$ws.Cells.Item(25, 9).Value = 0.8464338436636609
$ws.Cells.Item(25, 10).Value = 0.1294344242990055
$ws.Cells.Item(25, 12).Value = 0.3766005252823987
$ws.Cells.Item(25, 15).Value = 3.670893321681689
